$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Extend column BT and BU by copying formatting from column BS
#    (skip row 12, which has no data beyond column AX)
$ws.Range("BS1:BS11").Copy($ws.Range("BT1:BT11"))
$ws.Range("BS1:BS11").Copy($ws.Range("BU1:BU11"))
$ws.Range("BS13:BS29").Copy($ws.Range("BT13:BT29"))
$ws.Range("BS13:BS29").Copy($ws.Range("BU13:BU29"))

# 2. New header dates (2025-10-30 and 2025-10-31)
$ws.Range("BT1").Value = 45960
$ws.Range("BU1").Value = 45961

# 3. New attendance values per player row
$ws.Range("BT2").Value = "P"
$ws.Range("BU2").Value = "P"
$ws.Range("BT3").Value = "P"
$ws.Range("BU3").Value = "R"
$ws.Range("BT4").Value = "P"
$ws.Range("BU4").Value = "P"
$ws.Range("BT5").Value = "REP"
$ws.Range("BU5").Value = "P"
$ws.Range("BT6").Value = "B"
$ws.Range("BU6").Value = "B"
$ws.Range("BT7").Value = "P"
$ws.Range("BU7").Value = "REP"
$ws.Range("BT8").Value = "B"
$ws.Range("BU8").Value = "B"
$ws.Range("BT9").Value = "P"
$ws.Range("BU9").Value = "P"
$ws.Range("BT10").Value = "P"
$ws.Range("BU10").Value = "P"
$ws.Range("BT11").Value = "P"
$ws.Range("BU11").Value = "P"
$ws.Range("BT13").Value = "B"
$ws.Range("BU13").Value = "B"
$ws.Range("BT14").Value = "P"
$ws.Range("BU14").Value = "P"
$ws.Range("BT15").Value = "P"
$ws.Range("BU15").Value = "P"
$ws.Range("BT16").Value = "P"
$ws.Range("BU16").Value = "P"
$ws.Range("BT17").Value = "P"
$ws.Range("BU17").Value = "P"
$ws.Range("BT18").Value = "P"
$ws.Range("BU18").Value = "P"
$ws.Range("BT19").Value = "A"
$ws.Range("BU19").Value = "P"
$ws.Range("BT20").Value = "P"
$ws.Range("BU20").Value = "P"
$ws.Range("BT21").Value = "B"
$ws.Range("BU21").Value = "B"
$ws.Range("BT22").Value = "REP"
$ws.Range("BU22").Value = "P"
$ws.Range("BT23").Value = "RH"
$ws.Range("BU23").Value = "RH"
$ws.Range("BT24").Value = "P"
$ws.Range("BU24").Value = "P"
$ws.Range("BT25").Value = "A"
$ws.Range("BU25").Value = "A"
$ws.Range("BT26").Value = "P"
$ws.Range("BU26").Value = "M"
$ws.Range("BT27").Value = "REP"
$ws.Range("BU27").Value = "P"
$ws.Range("BT28").Value = "P"
$ws.Range("BU28").Value = "P"
$ws.Range("BT29").Value = "P"
$ws.Range("BU29").Value = "B"

# 4. Update the active selection to match the edited workbook state
$ws.Range("BW22").Select()
$excel.Calculate()
